# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-record rows (2-20) of the sheet:
# each destination row ends up with the full data record (Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg/unidad) that used to belong to a different row.
# Row 2 and row 15 keep their own record; the rest are permuted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that make up one data record (besides the identifying columns
# A,B,C,E,F,G,H,I,J,K which stay constant for every row in this sheet).
$cols = @("D","L","M","N","O","P","Q","R","S","T")

# Snapshot the "before" values for every row (2-20) and every tracked
# column so that later writes don't clobber data still needed as a
# source for another row.
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (the record that should end up
# living in the destination row), derived from the diff.
$map = @{
    2  = 2
    3  = 6
    4  = 5
    5  = 10
    6  = 17
    7  = 18
    8  = 19
    9  = 8
    10 = 4
    11 = 12
    12 = 16
    13 = 3
    14 = 11
    15 = 15
    16 = 7
    17 = 9
    18 = 13
    19 = 20
    20 = 14
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
